# Apply the data refresh to the cryptos sheet (rows 2-51), matching the
# upstream GitHub Actions scrape that produced the new Price/Volume(1h) values
# and re-ordered a couple of rows (Dai/RenderToken, Arweave/OKB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.691.60'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.792.31'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.03'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.80'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').Value = '3.790.74'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.34'
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.05'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '4.427.90'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '3.766.42'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.62'
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').Value = '67.721.15'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.07'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').Value = '  -9.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '459.87'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.700'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.34'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.02'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('E27').Value = '  -3.69%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '3.939.73'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.62'
$ws.Range('E34').Value = '  -2.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.07'
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.100'
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.78'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.15'
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.88'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.297'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '150.50'
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.29'
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.83'
$ws.Range('E49').Value = '  +4.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '389.37'
$ws.Range('E50').Value = '  -1.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.82'
$ws.Range('E51').Value = '  -5.14%  '
